$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Actividad"
$ws.Range("B1").Value = "Tipo de Consumo"
$ws.Range("C1").Value = "Unidad"
$ws.Range("D1").Value = "Alcance"
$ws.Range("E1").Value = "Valor"
$ws.Range("F1").Value = "Periodicidad"
$ws.Range("G1").Value = "Periodo de imputacion"

# Row 2
$ws.Range("A2").Value = "Combustion fija"
$ws.Range("B2").Value = "Gas Natural"
$ws.Range("C2").Value = "m3"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 20
$ws.Range("F2").Value = "Mensual"
$ws.Range("G2").Value = "02/2020"

# Row 3
$ws.Range("A3").Value = "Combustion movil"
$ws.Range("B3").Value = "Combustible consumido - Gasoil"
$ws.Range("C3").Value = "lts"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 40
$ws.Range("F3").Value = "Anual"
$ws.Range("G3").Value = "/2021"

# Row 4
$ws.Range("A4").Value = "Electricidad adquirida y consumida"
$ws.Range("B4").Value = "Electricidad"
$ws.Range("C4").Value = "Kwh"
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 200
$ws.Range("F4").Value = "Anual"
$ws.Range("G4").Value = "/2421"

# Row 5
$ws.Range("A5").Value = "Logistica de productos y residuos"
$ws.Range("B5").Value = "Camion de carga"
$ws.Range("C5").Value = "-"
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 40
$ws.Range("F5").Value = "Anual"
$ws.Range("G5").Value = "/2021"

# Row 6
$ws.Range("A6").Value = "Combustion fija"
$ws.Range("B6").Value = "Nafta"
$ws.Range("C6").Value = "lt"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 31
$ws.Range("F6").Value = "Mensual"
$ws.Range("G6").Value = "02/2026"

# Row 7
$ws.Range("A7").Value = "Logistica de productos y residuos"
$ws.Range("B7").Value = "Materia prima"
$ws.Range("C7").Value = "-"
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 100
$ws.Range("F7").Value = "Anual"
$ws.Range("G7").Value = "/2021"

# Remove any leftover rows below row 7 from the previous longer sheet (rows 8-14)
$ws.Range("A8:G20").Clear()

# Column G formatted as text ("@") so period values are not reinterpreted
$ws.Range("G1:G7").NumberFormat = "@"

# Wrap text for the long "Electricidad adquirida y consumida" activity cell
$ws.Range("A4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 45

# Column widths (bestFit sizing approximation)
$ws.Columns.Item(1).ColumnWidth = 29.6666666666667
$ws.Columns.Item(2).ColumnWidth = 29
$ws.Columns.Item(6).ColumnWidth = 11.3333333333333
$ws.Columns.Item(7).ColumnWidth = 20.6666666666667

# Restore selection seen in target workbook
$ws.Range("I11").Select()
